# test_server.xlsx edit: add a new "can match body" scenario row, update the
# headers-token fixture text, widen the query column, and flip which sheet /
# cell is the active selection when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: insert a new table row (row 5) for the "can match body" case ---
$ws1.Rows.Item(5).Insert()

$ws1.Range("A5").Value = "can match body"
$ws1.Range("B5").Value = "GET"
$ws1.Range("C5").Value = "/users/login"
$ws1.Range("E5").Value = "{`n  ""username"": ""user.maker"",`n  ""password"": ""1235""`n}"
$ws1.Range("G5").Value = 200
$ws1.Range("H5").Value = '{ "confirmationMessage": "Authorized user"}'

# Keep the table (ListObject) range / autofilter in sync with the new row.
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A1:H13"))

# --- Sheet1: update the "can match headers" fixture token value ---
$ws1.Range("D12").Value = '{"Authorization": "Token TheToken"}'

# --- Sheet1: widen the query column and move the active selection ---
$ws1.Columns.Item(6).ColumnWidth = 31.28515625
$ws1.Range("D13").Select()

# --- Make Sheet1 the active/selected tab instead of Sheet2 ---
$ws1.Activate()

# --- Sheet2: selection stays put (F18); it is simply no longer the active tab ---
$ws2.Range("F18").Select()
$ws1.Activate()
